$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Neurology" to "Session"
$ws.Name = "Session"

# The last logged row (row 42: student 212024 / Scan / 11:12:55) was removed.
# Deleting the entire row shifts everything below it up (there is nothing
# below it here) and shrinks the sheet's used range from A1:F42 to A1:F41.
$ws.Rows.Item(42).Delete()
